# Atualizei dados da bibi e add
# Insere uma nova linha de dados (dia 17 de Julho/2025) na planilha de
# faturamento diário, empurrando as linhas existentes (a partir da antiga
# linha 18) uma posição para baixo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere uma nova linha inteira na posição 18 (logo após o último dia de
# julho existente), deslocando para baixo tudo que estava nas linhas 18+.
$ws.Rows.Item(18).Insert()

# Preenche a nova linha com o novo registro de faturamento diário.
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 17617.07
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 2025
$ws.Range("E18").Value = "07/2025"
